$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 316; existing rows 316:410 shift down to 317:411.
$ws.Rows("316:316").Insert()

# Populate the newly inserted row 316 with the new data record.
$ws.Range("A316").Value = 4
$ws.Range("B316").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C316").Value = "Los Lagos"
$ws.Range("D316").Value = 45093
$ws.Range("E316").Value = 10
$ws.Range("F316").Value = 100112032
$ws.Range("G316").Value = "Zapallo italiano"
$ws.Range("H316").Value = "Sin especificar"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 240
$ws.Range("K316").Value = 17000
$ws.Range("L316").Value = 17000
$ws.Range("M316").Value = 17000
$ws.Range("N316").Value = '$/caja 50 unidades'
$ws.Range("O316").Value = "Región de Arica y Parinacota"
$ws.Range("P316").Value = 340
$ws.Range("Q316").Value = 50
$ws.Range("R316").Value = "Hortaliza"
